$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.555.44"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "1.846.06"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -1.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.67"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4641"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3854"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.83"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07892"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9967"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.49"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.857.83"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.948"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.123"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.75"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06659"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.11"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "27.553.83"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.387"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "2.072.91"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.80"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.49"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.113"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.408"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.77"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9755"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09403"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.589"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.298"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06034"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.315"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.180"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5888"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1867"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.35"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.243"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5591"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.909"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06698"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.99"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.051"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -1.27%  "
